{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n\nconst body = context.document.body;\n\n// ---------------------------------------------------------------------\n// 1) Update the first body paragraph (\"I am writing with the intention...\")\n//    so it ends with the new \"wondering if you would be interested...\"\n//    sentence instead of the old \"hopeful that you might consider...\" one.\n// ---------------------------------------------------------------------\nconst oldPara1 =\n  \"I am writing with the intention of applying to the University of California, \" +\n  \"Berkeley as a visiting student researcher. My research project aims to infer \" +\n  \"metabolic momentum based on moment differences in mass-weighted intensity \" +\n  \"distributions. I am hopeful that you might consider a collaborative venture \" +\n  \"with me in this regard.\";\nconst newPara1 =\n  \"I am writing with the intention of applying to the University of California, \" +\n  \"Berkeley as a visiting student researcher. My research project aims to infer \" +\n  \"metabolic momentum based on moment differences in mass-weighted intensity \" +\n  \"distributions. I am wondering if you would be interested in collaborating \" +\n  \"with me on my research project, and willing to serve as my faculty sponsor \" +\n  \"while I am there.\";\n\nconst hits1 = body.search(oldPara1, { matchCase: true });\nhits1.load(\"items\");\nawait context.sync();\nif (hits1.items.length > 0) {\n  hits1.items[0].insertText(newPara1, \"Replace\");\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 2) Replace the \"Additionally, it would be of immense benefit...\" paragraph\n//    with the (slightly shortened, \"completed\" removed) \"Enclosed, please\n//    find...\" paragraph text.\n// ---------------------------------------------------------------------\nconst oldPara2 =\n  \"Additionally, it would be of immense benefit to my application and \" +\n  \"subsequent research endeavors if you could agree to act as my faculty \" +\n  \"sponsor at Berkeley.\";\nconst newPara2 =\n  \"Enclosed, please find the concept note of my research project, my passport \" +\n  \"photograph, proof of English proficiency, a Visiting Student Researcher \" +\n  \"Appointment Request, and a Guarantee of Financial Support document for \" +\n  \"your perusal and convenience.\";\n\nconst hits2 = body.search(oldPara2, { matchCase: true });\nhits2.load(\"items\");\nawait context.sync();\nif (hits2.items.length > 0) {\n  hits2.items[0].insertText(newPara2, \"Replace\");\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 3) Replace the original \"Enclosed, please find...\" paragraph (the one\n//    that still says \"a completed Visiting Student Researcher Appointment\n//    Request\") with the \"Should you require any additional information...\"\n//    paragraph text.\n// ---------------------------------------------------------------------\nconst oldPara3 =\n  \"Enclosed, please find the concept note of my research project, my passport \" +\n  \"photograph, proof of English proficiency, a completed Visiting Student \" +\n  \"Researcher Appointment Request, and a Guarantee of Financial Support \" +\n  \"document for your perusal and convenience.\";\nconst newPara3 =\n  \"Should you require any additional information or have any questions, I am \" +\n  \"readily available to provide it. Your guidance and support would be \" +\n  \"invaluable to me, and I eagerly await the possibility of working under \" +\n  \"your esteemed mentorship.\";\n\nconst hits3 = body.search(oldPara3, { matchCase: true });\nhits3.load(\"items\");\nawait context.sync();\nif (hits3.items.length > 0) {\n  hits3.items[0].insertText(newPara3, \"Replace\");\n  await context.sync();\n}\n\n// ---------------------------------------------------------------------\n// 4) Remove the now-duplicated \"Should you require any additional\n//    information...\" paragraph (the original trailing one, which used to\n//    follow the blank paragraph carrying the \"_GoBack\" bookmark) together\n//    with the blank paragraph that immediately followed it.\n// ---------------------------------------------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nlet targetIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === newPara3) {\n    // The paragraph right after the \"_GoBack\" bookmark blank paragraph is\n    // the SECOND occurrence of this sentence; the first is the one we just\n    // wrote in step 3 and must be kept.\n    if (targetIndex === -1) {\n      targetIndex = -2; // mark \"first occurrence seen, skip it\"\n    } else if (targetIndex === -2) {\n      targetIndex = i;\n      break;\n    }\n  }\n}\n\nif (targetIndex >= 0) {\n  // Delete the blank paragraph that follows the duplicate first (so the\n  // duplicate paragraph's index stays valid), then delete the duplicate\n  // paragraph itself.\n  if (targetIndex + 1 < items.length && items[targetIndex + 1].text === \"\") {\n    items[targetIndex + 1].delete();\n  }\n  items[targetIndex].delete();\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the live document.\n\n$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------\n# 1) Remove the trailing duplicate-to-be \"Should you require any additional\n#    information...\" paragraph (originally paragraph #9, 1-indexed) and the\n#    blank paragraph right after it (originally paragraph #10). These are\n#    the last two paragraphs of the block that gets folded away once every\n#    paragraph below shifts up one slot. Delete the higher index first so\n#    the lower index stays valid.\n# ---------------------------------------------------------------------\n$d.Paragraphs(10).Range.Delete()\n$d.Paragraphs(9).Range.Delete()\n\n# ---------------------------------------------------------------------\n# 2) Update the \"I am writing with the intention of applying...\" paragraph\n#    so it ends with the new \"wondering if you would be interested...\"\n#    sentence instead of the old \"hopeful that you might consider...\" one.\n# ---------------------------------------------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"I am writing with the intention of applying to the University of California, Berkeley as a visiting student researcher. My research project aims to infer metabolic momentum based on moment differences in mass-weighted intensity distributions. I am hopeful that you might consider a collaborative venture with me in this regard.\"\n$find.Replacement.Text = \"I am writing with the intention of applying to the University of California, Berkeley as a visiting student researcher. My research project aims to infer metabolic momentum based on moment differences in mass-weighted intensity distributions. I am wondering if you would be interested in collaborating with me on my research project, and willing to serve as my faculty sponsor while I am there.\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# ---------------------------------------------------------------------\n# 3) Replace the \"Additionally, it would be of immense benefit...\" paragraph\n#    with the (slightly shortened, \"completed\" removed) \"Enclosed, please\n#    find...\" paragraph text.\n# ---------------------------------------------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Additionally, it would be of immense benefit to my application and subsequent research endeavors if you could agree to act as my faculty sponsor at Berkeley.\"\n$find.Replacement.Text = \"Enclosed, please find the concept note of my research project, my passport photograph, proof of English proficiency, a Visiting Student Researcher Appointment Request, and a Guarantee of Financial Support document for your perusal and convenience.\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n# ---------------------------------------------------------------------\n# 4) Replace the original \"Enclosed, please find...\" paragraph (the one\n#    that still says \"a completed Visiting Student Researcher Appointment\n#    Request\") with the \"Should you require any additional information...\"\n#    paragraph text.\n# ---------------------------------------------------------------------\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Enclosed, please find the concept note of my research project, my passport photograph, proof of English proficiency, a completed Visiting Student Researcher Appointment Request, and a Guarantee of Financial Support document for your perusal and convenience.\"\n$find.Replacement.Text = \"Should you require any additional information or have any questions, I am readily available to provide it. Your guidance and support would be invaluable to me, and I eagerly await the possibility of working under your esteemed mentorship.\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n"}
